$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the data between row 2 and row 5 (Date, Volumen, Precio minimo,
# Precio maximo, Precio promedio ponderado, Precio $/Kg)

# --- Save row 2's current values ---
$d2 = $ws.Range("D2").Value2
$m2 = $ws.Range("M2").Value2
$n2 = $ws.Range("N2").Value2
$o2 = $ws.Range("O2").Value2
$p2 = $ws.Range("P2").Value2
$s2 = $ws.Range("S2").Value2

# --- Save row 5's current values ---
$d5 = $ws.Range("D5").Value2
$m5 = $ws.Range("M5").Value2
$n5 = $ws.Range("N5").Value2
$o5 = $ws.Range("O5").Value2
$p5 = $ws.Range("P5").Value2
$s5 = $ws.Range("S5").Value2

# --- Write row 5's old values into row 2 ---
$ws.Range("D2").Value = $d5
$ws.Range("M2").Value = $m5
$ws.Range("N2").Value = $n5
$ws.Range("O2").Value = $o5
$ws.Range("P2").Value = $p5
$ws.Range("S2").Value = $s5

# --- Write row 2's old values into row 5 ---
$ws.Range("D5").Value = $d2
$ws.Range("M5").Value = $m2
$ws.Range("N5").Value = $n2
$ws.Range("O5").Value = $o2
$ws.Range("P5").Value = $p2
$ws.Range("S5").Value = $s2
